$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force text interpretation (leading apostrophe, like typing
    # into the Excel UI), then clear the implicit "Text" number
    # format Excel applies so the cell keeps the workbook default
    # style, matching the original (style-less) inline-string cell.
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "62.184.09"
$ws.Range("E2").Value = "  +2.96%  "

Set-TextCell "D3" "3.415.81"
$ws.Range("E3").Value = "  +3.30%  "

Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "

Set-TextCell "D5" "578.58"
$ws.Range("E5").Value = "  +2.80%  "

Set-TextCell "D6" "139.00"
$ws.Range("E6").Value = "  +7.07%  "

Set-TextCell "D7" "1.00"
$ws.Range("E7").Value = "  +0.03%  "

Set-TextCell "D8" "3.413.59"
$ws.Range("E8").Value = "  +3.10%  "

Set-TextCell "D9" "0.479"
$ws.Range("E9").Value = "  +1.40%  "

Set-TextCell "D10" "7.51"
$ws.Range("E10").Value = "  +0.88%  "

$ws.Range("E11").Value = "  +9.60%  "

$ws.Range("E12").Value = "  +6.37%  "

Set-TextCell "D13" "3.999.69"
$ws.Range("E13").Value = "  +3.54%  "

$ws.Range("E14").Value = "  +1.94%  "

Set-TextCell "D15" "0.0000182"
$ws.Range("E15").Value = "  +8.84%  "

Set-TextCell "D16" "3.416.09"
$ws.Range("E16").Value = "  +3.53%  "

Set-TextCell "D17" "25.55"
$ws.Range("E17").Value = "  +5.45%  "

Set-TextCell "D18" "62.220.71"
$ws.Range("E18").Value = "  +2.70%  "

Set-TextCell "D19" "14.18"
$ws.Range("E19").Value = "  +6.64%  "

$ws.Range("E20").Value = "  +4.31%  "

Set-TextCell "D21" "9.51"
$ws.Range("E21").Value = "  +6.25%  "

Set-TextCell "D22" "392.91"
$ws.Range("E22").Value = "  +11.95%  "

$ws.Range("E23").Value = "  +3.68%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell "D24" "3.554.36"
$ws.Range("E24").Value = "  +3.55%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D25" "0.0000129"
$ws.Range("E25").Value = "  +19.82%  "

$ws.Range("E26").Value = "  -0.01%  "

Set-TextCell "D27" "71.90"
$ws.Range("E27").Value = "  +3.95%  "

Set-TextCell "D28" "1.60"
$ws.Range("E28").Value = "  +9.85%  "

$ws.Range("E29").Value = "  +3.81%  "

Set-TextCell "D30" "0.999"
$ws.Range("E30").Value = "  +0.03%  "

Set-TextCell "D31" "8.35"
$ws.Range("E31").Value = "  +6.71%  "

Set-TextCell "D32" "0.161"
$ws.Range("E32").Value = "  +5.62%  "

$ws.Range("E33").Value = "  +3.54%  "

Set-TextCell "D34" "3.447.96"
$ws.Range("E34").Value = "  +3.47%  "

$ws.Range("E35").Value = "  +0.02%  "

Set-TextCell "D36" "23.65"
$ws.Range("E36").Value = "  +4.50%  "

Set-TextCell "D37" "5.53"
$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("E38").Value = "  +3.33%  "

$ws.Range("E39").Value = "  +6.25%  "

Set-TextCell "D40" "161.98"
$ws.Range("E40").Value = "  +3.06%  "

Set-TextCell "D41" "0.0796"
$ws.Range("E41").Value = "  +5.80%  "

$ws.Range("E42").Value = "  +14.75%  "

$ws.Range("E43").Value = "  +6.60%  "

$ws.Range("E44").Value = "  +0.03%  "

Set-TextCell "D45" "0.780"
$ws.Range("E45").Value = "  +4.92%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D46" "25.40"
$ws.Range("E46").Value = "  +12.46%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D47" "4.50"
$ws.Range("E47").Value = "  +3.79%  "

Set-TextCell "D48" "41.72"
$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("E49").Value = "  +4.72%  "

Set-TextCell "D50" "23.13"
$ws.Range("E50").Value = "  +6.19%  "

Set-TextCell "D51" "2.408.53"
$ws.Range("E51").Value = "  +11.90%  "
